# Add 2022-Q3 data
#
# Before: 总计, 2022-Q2, 2022-Q1
# After:  总计, 2022-Q3, 2022-Q2, 2022-Q1   (a new most-recent quarter is
#         inserted; the former "2022-Q2" sheet becomes the new "2022-Q3"
#         sheet with refreshed figures, and a brand-new "2022-Q2" sheet is
#         created holding the data the old "2022-Q2" sheet used to have)

$wb = $excel.ActiveWorkbook

# Assign a numeric-looking string to a cell while forcing Excel to keep
# it stored as TEXT (matching the source data, which keeps figures such
# as "5.54" as strings, not numbers) and without leaving the cell's
# number format / style changed.
function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: push the two existing quarters down one row
#    and put the new 2022-Q3 figures on top.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Make room for row 4 (copy row 3's "A" cell formatting onto the new row).
$summary.Range("A3").Copy()
$summary.Range("A4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Row 4 <- what used to be row 3 (2022-Q1, unchanged).
$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 1
$summary.Range("D4").Value = 0.25

# Row 3 <- what used to be row 2 (2022-Q2, unchanged).
$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.22

# Row 2 <- brand-new 2022-Q3 figures.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.21

# ---------------------------------------------------------------------
# 2) Quarter detail sheets.
# ---------------------------------------------------------------------
$oldQ2 = $wb.Worksheets.Item(2)   # currently named "2022-Q2"

# Duplicate the current "2022-Q2" sheet (with its untouched data) so the
# duplicate can become the new "2022-Q2" sheet once the original is
# turned into "2022-Q3".
$oldQ2.Copy([System.Reflection.Missing]::Value, $oldQ2)
$newQ2 = $wb.Worksheets.Item(3)

# Rename the original sheet first, otherwise the rename below collides
# with the still-existing "2022-Q2" name.
$oldQ2.Name = "2022-Q3"
$newQ2.Name = "2022-Q2"

# Refresh the (now renamed) sheet's fund figures.
$oldQ2.Range("C2").Value = "华安国际龙头（DAX）ETF（QDII）"
Set-TextValue $oldQ2.Range("D2") "5.54"
Set-TextValue $oldQ2.Range("E2") "93.57"
Set-TextValue $oldQ2.Range("F2") "3.80"
Set-TextValue $oldQ2.Range("G2") "0.2105"

# Restore "2022-Q1" (unchanged from before the edit) as the active tab,
# matching its original selected state.
$wb.Worksheets.Item("2022-Q1").Activate()

Write-Output "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Output $s.Name
}
